# Daily attendance processing - 2025-11-12 05:23:51
#
# Normalize the "Recorded By" audit trail in column G of the
# "Session Analysis Results" sheet: move the "System" / "system" token
# that currently leads each comma-separated list to the end of that
# list instead.
#
#   "System, <name>"                  -> "<name>, System"
#   "System, system, <name>"          -> "System, <name>, system"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows whose Recorded By value is exactly "System, dnasr281@gmail.com"
# -> reorder to "dnasr281@gmail.com, System"
$rowsTwoPart = @(3, 6, 10, 12, 13, 14, 15, 18, 19, 20, 21, 22, 24, 26, 29, 32, 36, 38, 39, 40, 41, 44, 45, 46, 47, 48, 50, 52, 55, 58, 62, 64, 65, 66, 67, 70, 71, 72, 73, 74, 76, 78, 83, 84, 85, 86, 90, 92, 99, 101, 109, 110, 111, 112, 116, 118, 125, 127, 135, 136, 137, 138, 142, 144, 151, 153)

foreach ($r in $rowsTwoPart) {
    $ws.Range("G$r").Value = "dnasr281@gmail.com, System"
}

# Rows whose Recorded By value is exactly "System, system, backup@backdoor.com"
# -> reorder to "System, backup@backdoor.com, system"
$rowsThreePart = @(2, 28, 54)

foreach ($r in $rowsThreePart) {
    $ws.Range("G$r").Value = "System, backup@backdoor.com, system"
}
